$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Natmi LR-pair table (Il6 -> Il6st) was recomputed with an extra replicate/
# sample pooled in ("Dr Hou advice"), turning the old 4 sending-clusters x
# 3 target-clusters (12 data rows) layout into a full 4x4 cross of clusters
# (ECs, FAPs, M2, sCs) -> 16 data rows, rows 2-17.

$rows = @(
    @("ECs", "Il6", "Il6st", "ECs", 2, 0.6666666666666666, 11.898632, 35.695896, 0.1819741108242693, 0.1819741108242693, 3, 1, 44.017783, 132.053349, 0.253483683026081, 0.253483683026081, 523.7514013728561, 4713.762612355705, 0.04612746782713202, 0.04612746782713201),
    @("ECs", "Il6", "Il6st", "FAPs", 2, 0.6666666666666666, 11.898632, 35.695896, 0.1819741108242693, 0.1819741108242693, 3, 1, 99.76728333333334, 299.30185, 0.5745264004968147, 0.5745264004968147, 1187.094190023067, 10683.8477102076, 0.1045489308754759, 0.1045489308754759),
    @("ECs", "Il6", "Il6st", "M2", 2, 0.6666666666666666, 11.898632, 35.695896, 0.1819741108242693, 0.1819741108242693, 3, 1, 10.852095, 32.556285, 0.06249358376701795, 0.06249358376701795, 129.12508483404, 1162.12576350636, 0.01137221433822508, 0.01137221433822508),
    @("ECs", "Il6", "Il6st", "sCs", 2, 0.6666666666666666, 11.898632, 35.695896, 0.1819741108242693, 0.1819741108242693, 3, 1, 19.01418566666667, 57.042557, 0.1094963327100864, 0.1094963327100864, 226.2427980273413, 2036.185182246072, 0.01992549778343633, 0.01992549778343633),
    @("FAPs", "Il6", "Il6st", "ECs", 3, 1, 50.04372666666666, 150.13118, 0.7653537534818657, 0.7653537534818657, 3, 1, 44.017783, 132.053349, 0.253483683026081, 0.253483683026081, 2202.813900924647, 19825.32510832182, 0.1940046882504186, 0.1940046882504186),
    @("FAPs", "Il6", "Il6st", "FAPs", 3, 1, 50.04372666666666, 150.13118, 0.7653537534818657, 0.7653537534818657, 3, 1, 99.76728333333334, 299.30185, 0.5745264004968147, 0.5745264004968147, 4992.726657409222, 44934.539916683, 0.4397159370946628, 0.4397159370946628),
    @("FAPs", "Il6", "Il6st", "M2", 3, 1, 50.04372666666666, 150.13118, 0.7653537534818657, 0.7653537534818657, 3, 1, 10.852095, 32.556285, 0.06249358376701795, 0.06249358376701795, 543.0792759407, 4887.7134834663, 0.04782969890462058, 0.04782969890462058),
    @("FAPs", "Il6", "Il6st", "sCs", 3, 1, 50.04372666666666, 150.13118, 0.7653537534818657, 0.7653537534818657, 3, 1, 19.01418566666667, 57.042557, 0.1094963327100864, 0.1094963327100864, 951.5407102919177, 8563.86639262726, 0.08380342923216383, 0.08380342923216383),
    @("M2", "Il6", "Il6st", "ECs", 3, 1, 0.8079146666666667, 2.423744, 0.01235600471453799, 0.01235600471453799, 3, 1, 44.017783, 132.053349, 0.253483683026081, 0.253483683026081, 35.56261247985067, 320.063512318656, 0.00313204558252871, 0.00313204558252871),
    @("M2", "Il6", "Il6st", "FAPs", 3, 1, 0.8079146666666667, 2.423744, 0.01235600471453799, 0.01235600471453799, 3, 1, 99.76728333333334, 299.30185, 0.5745264004968147, 0.5745264004968147, 80.60345145848889, 725.4310631264001, 0.007098850913165183, 0.007098850913165183),
    @("M2", "Il6", "Il6st", "M2", 3, 1, 0.8079146666666667, 2.423744, 0.01235600471453799, 0.01235600471453799, 3, 1, 10.852095, 32.556285, 0.06249358376701795, 0.06249358376701795, 8.767566714560001, 78.90810043104001, 0.0007721710156536484, 0.0007721710156536485),
    @("M2", "Il6", "Il6st", "sCs", 3, 1, 0.8079146666666667, 2.423744, 0.01235600471453799, 0.01235600471453799, 3, 1, 19.01418566666667, 57.042557, 0.1094963327100864, 0.1094963327100864, 15.36183947482311, 138.256555273408, 0.001352937203190448, 0.001352937203190448),
    @("sCs", "Il6", "Il6st", "ECs", 3, 1, 2.636126666666667, 7.90838, 0.040316130979327, 0.040316130979327, 3, 1, 44.017783, 132.053349, 0.253483683026081, 0.253483683026081, 116.0364515738467, 1044.32806416462, 0.01021948136600169, 0.01021948136600169),
    @("sCs", "Il6", "Il6st", "FAPs", 3, 1, 2.636126666666667, 7.90838, 0.040316130979327, 0.040316130979327, 3, 1, 99.76728333333334, 299.30185, 0.5745264004968147, 0.5745264004968147, 262.9991960558889, 2366.992764503, 0.02316268161351086, 0.02316268161351086),
    @("sCs", "Il6", "Il6st", "M2", 3, 1, 2.636126666666667, 7.90838, 0.040316130979327, 0.040316130979327, 3, 1, 10.852095, 32.556285, 0.06249358376701795, 0.06249358376701795, 28.6074970187, 257.4674731683, 0.002519499508518639, 0.00251949950851864),
    @("sCs", "Il6", "Il6st", "sCs", 3, 1, 2.636126666666667, 7.90838, 0.040316130979327, 0.040316130979327, 3, 1, 19.01418566666667, 57.042557, 0.1094963327100864, 0.1094963327100864, 50.12380188085111, 451.11421692766, 0.004414468491295812, 0.004414468491295812)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]
    $excelRow = $i + 2
    for ($c = 0; $c -lt $r.Length; $c++) {
        $ws.Cells.Item($excelRow, $c + 1).Value = $r[$c]
    }
}